$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values (cryptos.xlsx) scraped on Fri Feb 10 09:53:24 UTC 2023
$ws.Range("D2").Value = '''308.64'
$ws.Range("E2").Value = '''-4.15%'
$ws.Range("D3").Value = '''40.08'
$ws.Range("E3").Value = '''-6.03%'
$ws.Range("D4").Value = '''5.121'
$ws.Range("E4").Value = '''-0.74%'
$ws.Range("D5").Value = '''0.07741'
$ws.Range("E5").Value = '''-5.57%'
$ws.Range("D6").Value = '''4.257'
$ws.Range("E6").Value = '''-0.83%'
$ws.Range("D7").Value = '''1.600'
$ws.Range("E7").Value = '''-11.43%'
$ws.Range("E8").Value = '''-5.19%'
$ws.Range("D10").Value = '''0.1747'
$ws.Range("E10").Value = '''-6.58%'
$ws.Range("D11").Value = '''0.09038'
$ws.Range("E11").Value = '''-4.93%'
$ws.Range("D12").Value = '''0.04436'
$ws.Range("E12").Value = '''-5.00%'
$ws.Range("D13").Value = '''0.1055'
$ws.Range("E13").Value = '''-0.24%'
$ws.Range("D14").Value = '''0.001254'
$ws.Range("E14").Value = '''-2.70%'
$ws.Range("B15").Value = '''TigerCash'
$ws.Range("C15").Value = '''https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.005799'
$ws.Range("E15").Value = '''1.26%'
$ws.Range("B16").Value = '''UpBots'
$ws.Range("C16").Value = '''https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").Value = '''0.007491'
$ws.Range("E16").Value = '''2,416.67%'
$ws.Range("B17").Value = '''LEO'
$ws.Range("C17").Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.356'
$ws.Range("E17").Value = '''-0.24%'
$ws.Range("B18").Value = '''BTSEToken'
$ws.Range("C18").Value = '''https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.420'
$ws.Range("E18").Value = '''-4.28%'
$ws.Range("B19").Value = '''BitpandaEcosystemToken'
$ws.Range("C19").Value = '''https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '''0.3276'
$ws.Range("E19").Value = '''-2.93%'
$ws.Range("B20").Value = '''MCDex'
$ws.Range("C20").Value = '''https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D20").Value = '''7.083'
$ws.Range("E20").Value = '''-4.62%'
$ws.Range("B21").Value = '''ProBitToken'
$ws.Range("C21").Value = '''https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = '''0.1340'
$ws.Range("E21").Value = '''-3.53%'
$ws.Range("B22").Value = '''ZBToken'
$ws.Range("C22").Value = '''https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D22").Value = '''0.2788'
$ws.Range("E22").Value = '''11.88%'
$ws.Range("B23").Value = '''CoinExToken'
$ws.Range("C23").Value = '''https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").Value = '''0.04173'
$ws.Range("E23").Value = '''0.50%'
$ws.Range("D24").Value = '''0.001203'
$ws.Range("E24").Value = '''-3.40%'
$ws.Range("D25").Value = '''0.004092'
$ws.Range("E25").Value = '''-7.57%'
$ws.Range("D26").Value = '''0.0001304'
$ws.Range("E26").Value = '''8.72%'
$ws.Range("D38").Value = '''0.02366'
$ws.Range("E38").Value = '''-14.69%'
$ws.Range("D39").Value = '''0.05224'
$ws.Range("E39").Value = '''-6.64%'
$ws.Range("D40").Value = '''0.007943'
$ws.Range("E40").Value = '''-1.44%'
$ws.Range("D41").Value = '''0.1328'
$ws.Range("E41").Value = '''-4.99%'
$ws.Range("D42").Value = '''0.006493'
$ws.Range("E42").Value = '''-0.79%'
$ws.Range("D43").Value = '''0.001962'
$ws.Range("E43").Value = '''-5.82%'
$ws.Range("D44").Value = '''0.008775'
$ws.Range("E44").Value = '''5.39%'
$ws.Range("E45").Value = '''-4.85%'
$ws.Range("D46").Value = '''0.00006566'
$ws.Range("E46").Value = '''-5.55%'
$ws.Range("D47").Value = '''0.00000000752'
$ws.Range("E47").Value = '''0.35%'
$ws.Range("E48").Value = '''98.62%'
$ws.Range("D49").Value = '''0.004125'
$ws.Range("E49").Value = '''18.55%'
$ws.Range("D50").Value = '''0.00002106'
$ws.Range("E50").Value = '''0.35%'
$ws.Range("D51").Value = '''0.0002005'
$ws.Range("E51").Value = '''0.35%'
